$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the BOM table by Designator (column A) ascending so MISC1/MISC2 move
# to the top, along with their per-row formatting.
$sortRange = $ws.Range("A2:G9")
$sortRange.Sort($ws.Range("A2"))

# Add a bias resistor to the anode: extend the resistor designator list and
# bump the quantity from 6 to 12 (R1,R2,R5-R8 -> R1-R12).
$ws.Cells.Item(7,1).Value = "'R1, R2, R3, R4, R5, R6, R7, R8, R9, R10, R11, R12"
$ws.Cells.Item(7,6).Value = 12

# The MISC1/MISC2 rows carried a bold red font before the sort; normalize
# them to the same plain font used by every other data row.
$fmtSrc = $ws.Range("A4:G4")
$fmtSrc.Copy()
$ws.Range("A2:G2").PasteSpecial(-4122)
$fmtSrc.Copy()
$ws.Range("A3:G3").PasteSpecial(-4122)

# Match the narrower "Do Not Place" column width from the edit.
$ws.Columns.Item(7).ColumnWidth = 19.43

$ws.Range("A1").Select() | Out-Null
